$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.668100000000003
$ws.Range("A10").Value = -22.0494
$ws.Range("A12").Value = -21.41260000000001
$ws.Range("D12").Value = -5.873899999999995
$ws.Range("D17").Value = -8.277899999999988
$ws.Range("A18").Value = -22.07240000000001
$ws.Range("D26").Value = -7.393500000000006
$ws.Range("D27").Value = -8.1981
$ws.Range("D28").Value = -8.247999999999998
$ws.Range("A37").Value = -20.08259999999998
$ws.Range("D37").Value = -7.880200000000003
$ws.Range("A55").Value = -22.232
$ws.Range("D65").Value = -7.770000000000001
$ws.Range("A68").Value = -21.42150000000001
$ws.Range("D73").Value = -8.311799999999998
$ws.Range("A77").Value = -20.21919999999999
$ws.Range("A78").Value = -19.77139999999998
$ws.Range("D84").Value = -8.209700000000003
$ws.Range("D85").Value = -8.498799999999997
$ws.Range("D93").Value = -6.675099999999992
$ws.Range("D95").Value = -7.5172
$ws.Range("D98").Value = -7.403100000000006
$ws.Range("D99").Value = -8.030900000000003
$ws.Range("D101").Value = -8.168899999999999
